$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was added for this market/category subset.
# It belongs chronologically in the first data slot (row 107), so insert a
# new row there - this pushes the existing rows 107..231 down to 108..232
# (each row keeps its own full set of values, only its row number shifts).
$ws.Rows.Item(107).Insert()

# Populate the newly inserted row 107 with the new record.
$ws.Cells.Item(107, 1).Value  = 11
$ws.Cells.Item(107, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(107, 3).Value  = "Bíobío"
$ws.Cells.Item(107, 4).Value  = 44902
$ws.Cells.Item(107, 5).Value  = 8
$ws.Cells.Item(107, 6).Value  = 100112003
$ws.Cells.Item(107, 7).Value  = "Ajo"
$ws.Cells.Item(107, 8).Value  = "Chino"
$ws.Cells.Item(107, 9).Value  = "1a (cosecha)"
$ws.Cells.Item(107, 10).Value = 180
$ws.Cells.Item(107, 11).Value = 13500
$ws.Cells.Item(107, 12).Value = 14000
$ws.Cells.Item(107, 13).Value = 13722
$ws.Cells.Item(107, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(107, 15).Value = "China"
$ws.Cells.Item(107, 16).Value = 1372
$ws.Cells.Item(107, 17).Value = 10
$ws.Cells.Item(107, 18).Value = "Hortaliza"
